# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the formatting of the existing header cell (AC1) onto the three
# new header cells so they match the bold/centered/bordered header style.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the win/loss/tie record for every data row (2..45).
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 86   # AD = Wins
    $ws.Cells.Item($r, 31).Value = 76   # AE = Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF = Ties
}
